$wb = $excel.ActiveWorkbook

# ALC row 2 (G=5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179.8
$ws.Range("I2").Value = 179.8
$ws.Range("K2").Value = 179.8
$ws.Range("M2").Value = -66.80000000000001

# ALC row 18 (G=5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2513.7778
$ws.Range("I18").Value = 2690.5
$ws.Range("J18").Value = 1100
$ws.Range("K18").Value = 2690.5
$ws.Range("L18").Value = 1100
$ws.Range("M18").Value = -2406.5
$ws.Range("N18").Value = -1668

# ALC row 33 (G=5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 110.75
$ws.Range("I33").Value = 103.92857
$ws.Range("K33").Value = 103.92857
$ws.Range("M33").Value = 125.07143

# ALC row 40 (G=5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4968.7646
$ws.Range("I40").Value = 3407.7778
$ws.Range("J40").Value = 6724.875
$ws.Range("K40").Value = 3407.7778
$ws.Range("L40").Value = 6724.875
$ws.Range("M40").Value = -3232.7778
$ws.Range("N40").Value = -7074.875

# ALC row 101 (G=19884)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 3034.8235
$ws.Range("I101").Value = 517.8182
$ws.Range("J101").Value = 7649.3335
$ws.Range("K101").Value = 1553.4546
$ws.Range("L101").Value = 22948.0005
$ws.Range("M101").Value = 68.54539999999997
$ws.Range("N101").Value = -26192.0005

# ALC row 113 (G=27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 52634720
$ws.Range("I113").Value = 62502588
$ws.Range("K113").Value = 62502588
$ws.Range("M113").Value = -62499334

# ALC row 141 (G=44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4191.5
$ws.Range("I141").Value = 3849.182
$ws.Range("K141").Value = 11547.546
$ws.Range("M141").Value = -6367.545999999998

# ARM row 2 (G=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3355.3333
$ws.Range("I2").Value = 3355.3333
$ws.Range("K2").Value = 3355.3333
$ws.Range("M2").Value = -3242.3333

# ARM row 26 (G=2241)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1800
$ws.Range("I26").Value = 1625
$ws.Range("J26").Value = 2500
$ws.Range("K26").Value = 1625
$ws.Range("L26").Value = 2500
$ws.Range("M26").Value = -1295
$ws.Range("N26").Value = -3160

# ARM row 45 (G=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7308.7334
$ws.Range("I45").Value = 3272.8333
$ws.Range("K45").Value = 3272.8333
$ws.Range("M45").Value = -2895.8333

# ARM row 61 (G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35724550
$ws.Range("I61").Value = 62508500
$ws.Range("J61").Value = 12624.833
$ws.Range("K61").Value = 62508500
$ws.Range("L61").Value = 12624.833
$ws.Range("M61").Value = -62508288
$ws.Range("N61").Value = -13048.833

# ARM row 116 (G=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3355.3333
$ws.Range("I116").Value = 3355.3333
$ws.Range("K116").Value = 3355.3333
$ws.Range("M116").Value = -1061.3333

# ARM row 122 (G=36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1799.5
$ws.Range("I122").Value = 1799.5
$ws.Range("K122").Value = 5398.5
$ws.Range("M122").Value = -2948.5

# ARM row 136 (G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 35724550
$ws.Range("I136").Value = 62508500
$ws.Range("J136").Value = 12624.833
$ws.Range("K136").Value = 187525500
$ws.Range("L136").Value = 37874.499
$ws.Range("M136").Value = -187522950
$ws.Range("N136").Value = -42974.499

# BSM row 3 (G=27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3355.3333
$ws.Range("I3").Value = 3355.3333
$ws.Range("K3").Value = 3355.3333
$ws.Range("M3").Value = -3241.3333

# BSM row 76 (G=10630)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 5749.75
$ws.Range("J76").Value = 5749.75
$ws.Range("L76").Value = 5749.75
$ws.Range("N76").Value = -6379.75

# BSM row 79 (G=10630)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 5749.75
$ws.Range("J79").Value = 5749.75
$ws.Range("L79").Value = 5749.75
$ws.Range("N79").Value = -7933.75

# BSM row 80 (G=13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 632.25
$ws.Range("I80").Value = 395.75
$ws.Range("K80").Value = 395.75
$ws.Range("M80").Value = 602.25

# BSM row 83 (G=13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 632.25
$ws.Range("I83").Value = 395.75
$ws.Range("K83").Value = 1978.75
$ws.Range("M83").Value = 3013.25

# CRP row 31 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4410.1035
$ws.Range("I31").Value = 3699.6667
$ws.Range("J31").Value = 5171.2856
$ws.Range("K31").Value = 3699.6667
$ws.Range("L31").Value = 5171.2856
$ws.Range("M31").Value = -3404.6667
$ws.Range("N31").Value = -5761.2856

# CRP row 34 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4410.1035
$ws.Range("I34").Value = 3699.6667
$ws.Range("J34").Value = 5171.2856
$ws.Range("K34").Value = 3699.6667
$ws.Range("L34").Value = 5171.2856
$ws.Range("M34").Value = -3497.6667
$ws.Range("N34").Value = -5575.2856

# CRP row 62 (G=12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4332.6665
$ws.Range("I62").Value = 4499.5
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 4499.5
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -3875.5
$ws.Range("N62").Value = -5247

# CRP row 65 (G=12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4332.6665
$ws.Range("I65").Value = 4499.5
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 22497.5
$ws.Range("L65").Value = 19995
$ws.Range("M65").Value = -19377.5
$ws.Range("N65").Value = -26235

# CUL row 2 (G=4847)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9148.454
$ws.Range("I2").Value = 42
$ws.Range("K2").Value = 252
$ws.Range("M2").Value = -139

# CUL row 38 (G=4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 772.2727
$ws.Range("I38").Value = 18
$ws.Range("J38").Value = 1400.8334
$ws.Range("K38").Value = 54
$ws.Range("L38").Value = 4202.5002
$ws.Range("M38").Value = 293
$ws.Range("N38").Value = -4896.5002

# CUL row 122 (G=36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4928.6
$ws.Range("J122").Value = 1181
$ws.Range("L122").Value = 10629
$ws.Range("N122").Value = -15529

# GSM row 21 (G=4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# GSM row 30 (G=4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# LTW row 22 (G=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2250
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705

# LTW row 27 (G=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2250
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893

# LTW row 46 (G=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 16929.111
$ws.Range("J46").Value = 21209.572
$ws.Range("L46").Value = 21209.572
$ws.Range("N46").Value = -21585.572

# LTW row 136 (G=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6766.1113
$ws.Range("I136").Value = 6770.7144
$ws.Range("K136").Value = 20312.1432
$ws.Range("M136").Value = -17762.1432

# WVR row 81 (G=12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1659.5
$ws.Range("I81").Value = 1659.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3319
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2258
$ws.Range("N81").ClearContents()

# WVR row 84 (G=12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1659.5
$ws.Range("I84").Value = 1659.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 16595
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -11291
$ws.Range("N84").ClearContents()

# WVR row 107 (G=27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 702.8125
$ws.Range("I107").Value = 596
$ws.Range("J107").Value = 1279.6
$ws.Range("K107").Value = 1788
$ws.Range("L107").Value = 3838.8
$ws.Range("M107").Value = 132
$ws.Range("N107").Value = -7678.799999999999

# WVR row 122 (G=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2850.2727
$ws.Range("I122").Value = 2235.8096
$ws.Range("J122").Value = 3925.5833
$ws.Range("K122").Value = 6707.4288
$ws.Range("L122").Value = 11776.7499
$ws.Range("M122").Value = -4257.4288
$ws.Range("N122").Value = -16676.7499

# WVR row 132 (G=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3173.577
$ws.Range("I132").Value = 2835.5
$ws.Range("J132").Value = 4300.5
$ws.Range("K132").Value = 8506.5
$ws.Range("L132").Value = 12901.5
$ws.Range("M132").Value = -5976.5
$ws.Range("N132").Value = -17961.5
